$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.889.35'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '1.545.70'
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '205.80'
$ws.Range('E5').Value = '  -0.15%  '
$ws.Range('E6').Value = '  -0.65%  '
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.247'
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '21.27'
$ws.Range('E9').Value = '  -2.29%  '
$ws.Range('E10').Value = '  -0.42%  '
$ws.Range('E11').Value = '  -0.75%  '
$ws.Range('D12').Value = '1.764.76'
$ws.Range('E12').Value = '  -1.16%  '
$ws.Range('D13').Value = '1.544.05'
$ws.Range('E13').Value = '  -1.30%  '
$ws.Range('E15').Value = '  -0.78%  '
$ws.Range('D16').Value = '26.865.08'
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.41'
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '213.71'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.17'
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('E22').Value = '  -2.57%  '
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('E24').Value = '  -3.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.56'
$ws.Range('E25').Value = '  -0.66%  '
$ws.Range('E26').Value = '  -2.00%  '
$ws.Range('E27').Value = '  -0.87%  '
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0457'
$ws.Range('E30').Value = '  -1.75%  '
$ws.Range('E31').Value = '  -1.20%  '
$ws.Range('E32').Value = '  +1.41%  '
$ws.Range('D33').Value = '1.354.84'
$ws.Range('E33').Value = '  -3.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.94'
$ws.Range('E34').Value = '  +0.46%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.964'
$ws.Range('E36').Value = '  +4.74%  '
$ws.Range('E37').Value = '  +0.36%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.520'
$ws.Range('E39').Value = '  -1.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.804'
$ws.Range('E40').Value = '  -0.83%  '
$ws.Range('E41').Value = '  +0.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.57'
$ws.Range('E42').Value = '  +2.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.989'
$ws.Range('E43').Value = '  -0.88%  '
$ws.Range('E44').Value = '  +1.93%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '63.37'
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('E46').Value = '  -2.19%  '
$ws.Range('D47').Value = '1.678.83'
$ws.Range('E47').Value = '  -1.23%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.86'
$ws.Range('E48').Value = '  -0.58%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0510'
$ws.Range('E49').Value = '  +0.75%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₇0972'
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0946'
$ws.Range('E51').Value = '  -0.27%  '
